# Daily update at 8 AM UTC - append the next day's row of data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A59").Value = 46008
$ws.Range("A59").NumberFormat = $ws.Range("A58").NumberFormat
$ws.Range("B59").Value = 122
$ws.Range("C59").Value = 139
$ws.Range("D59").Value = 128
